$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '27.190.53'
$ws.Range('E2').Value = '  +0.28%  '
$ws.Range('D3').Value = '1.904.93'
$ws.Range('E3').Value = '  +0.88%  '
$ws.Range('E4').Value = '  +0.06%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '306.53'
$ws.Range('D5').ClearFormats()
$ws.Range('E6').Value = '  +0.02%  '
$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '0.5254'
$ws.Range('D7').ClearFormats()
$ws.Range('E7').Value = '  +2.02%  '
$ws.Range('E8').Value = '  +1.41%  '
$ws.Range('E9').Value = '  +0.58%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '21.15'
$ws.Range('D10').ClearFormats()
$ws.Range('E10').Value = '  +0.44%  '
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '0.9004'
$ws.Range('D11').ClearFormats()
$ws.Range('E11').Value = '  -0.30%  '
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '0.08405'
$ws.Range('D12').ClearFormats()
$ws.Range('E12').Value = '  +10.18%  '
$ws.Range('D13').Value = '1.906.40'
$ws.Range('E13').Value = '  +0.60%  '
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '95.05'
$ws.Range('D14').ClearFormats()
$ws.Range('E14').Value = '  +0.44%  '
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '5.284'
$ws.Range('D15').ClearFormats()
$ws.Range('E15').Value = '  +0.27%  '
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '0.000008613'
$ws.Range('D17').ClearFormats()
$ws.Range('E17').Value = '  +1.20%  '
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '14.56'
$ws.Range('D18').ClearFormats()
$ws.Range('E19').Value = '  +0.06%  '
$ws.Range('D20').Value = '27.233.33'
$ws.Range('E20').Value = '  +0.27%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '5.066'
$ws.Range('D21').ClearFormats()
$ws.Range('E21').Value = '  +0.31%  '
$ws.Range('D22').Value = '2.143.66'
$ws.Range('E22').Value = '  +0.91%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '10.62'
$ws.Range('D23').ClearFormats()
$ws.Range('E23').Value = '  +0.68%  '
$ws.Range('E24').Value = '  +0.19%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '147.56'
$ws.Range('D25').ClearFormats()
$ws.Range('E25').Value = '  +1.29%  '
$ws.Range('E26').Value = '  +4.91%  '
$ws.Range('E27').Value = '  -2.11%  '
$ws.Range('E28').Value = '  +0.70%  '
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '114.87'
$ws.Range('D29').ClearFormats()
$ws.Range('E29').Value = '  +0.19%  '
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '4.924'
$ws.Range('D30').ClearFormats()
$ws.Range('E30').Value = '  -1.25%  '
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '4.817'
$ws.Range('D31').ClearFormats()
$ws.Range('E31').Value = '  -0.16%  '
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '0.09284'
$ws.Range('D32').ClearFormats()
$ws.Range('E32').Value = '  +0.82%  '
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '0.8089'
$ws.Range('D33').ClearFormats()
$ws.Range('E33').Value = '  +6.04%  '
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '0.05066'
$ws.Range('D34').ClearFormats()
$ws.Range('E34').Value = '  +0.02%  '
$ws.Range('E35').Value = '  +3.52%  '
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '2.951'
$ws.Range('D36').ClearFormats()
$ws.Range('E36').Value = '  -1.68%  '
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '3.375'
$ws.Range('D37').ClearFormats()
$ws.Range('E37').Value = '  +2.94%  '
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '2.611'
$ws.Range('D38').ClearFormats()
$ws.Range('E38').Value = '  +1.43%  '
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '0.5739'
$ws.Range('D39').ClearFormats()
$ws.Range('E39').Value = '  +2.25%  '
$ws.Range('E40').Value = '  -0.24%  '
$ws.Range('E41').Value = '  -0.12%  '
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '6.654'
$ws.Range('D42').ClearFormats()
$ws.Range('E42').Value = '  +1.17%  '
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '8.985'
$ws.Range('D43').ClearFormats()
$ws.Range('E43').Value = '  -0.47%  '
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '117.49'
$ws.Range('D44').ClearFormats()
$ws.Range('E44').Value = '  -1.05%  '
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '0.1513'
$ws.Range('D45').ClearFormats()
$ws.Range('E45').Value = '  +0.38%  '
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '0.4854'
$ws.Range('D46').ClearFormats()
$ws.Range('E46').Value = '  +1.12%  '
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '10.17'
$ws.Range('D47').ClearFormats()
$ws.Range('E47').Value = '  -0.18%  '
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '1.0000'
$ws.Range('D48').ClearFormats()
$ws.Range('E48').Value = '  +0.01%  '
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '1.617'
$ws.Range('D49').ClearFormats()
$ws.Range('E49').Value = '  +2.41%  '
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '37.48'
$ws.Range('D50').ClearFormats()
$ws.Range('E50').Value = '  +0.67%  '
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '63.82'
$ws.Range('D51').ClearFormats()
$ws.Range('E51').Value = '  +0.21%  '
